$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# Update the daily conversion note text on Hoja1!A1
$old1 = [char]0x2705 + " 1000 Bs = 2.99 = 11056.87 pesos"
$new1 = [char]0x2705 + " 1000 Bs = 2.93 = 10856.61 pesos"
$old2 = [char]0x2705 + " 11056.87 pesos = 2.96 = 953.55 Bs"
$new2 = [char]0x2705 + " 10856.61 pesos = 2.93 = 955.32 Bs"

$text = $ws1.Range("A1").Text
$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$ws1.Range("A1").Value = $text

# Update rate cells on the "tasas" sheet
$ws2.Range("N10").Value = 340.995
$ws2.Range("O10").Value = 3702.05
$ws2.Range("N12").Value = 3699.09
$ws2.Range("O12").Value = 325.5
